$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated avg_potential_profit values for C2:C31
$values = @(
    1051105.8799999999,
    1043333.33,
    1012854.55,
    1018576.67,
    1034258.33,
    1044844.44,
    902371.43,
    1029000,
    1031263.64,
    995290.91,
    1060000,
    976900,
    995200,
    1038400,
    1033450,
    1040100,
    1060000,
    1060900,
    1027600,
    1060000,
    973600,
    1016800,
    962800,
    998500,
    844000,
    915400,
    1060000,
    952000,
    978600,
    1041800
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 3).Value = $val
    $row++
}

# Update selected cell to E19
$ws.Range("E19").Select()
